$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Rarres2"
$ws.Range("C2").Value = "Cmklr1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 1.445726333333333
$ws.Range("H2").Value = 4.337179
$ws.Range("I2").Value = 0.01544456920939864
$ws.Range("J2").Value = 0.01544456920939864
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.803353666666666
$ws.Range("N2").Value = 8.410060999999999
$ws.Range("O2").Value = 0.03356888849358693
$ws.Range("P2").Value = 0.03356888849358693
$ws.Range("Q2").Value = 4.052882217546554
$ws.Range("R2").Value = 36.475939957919
$ws.Range("S2").Value = 0.0005184570216217889
$ws.Range("T2").Value = 0.0005184570216217889
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Rarres2"
$ws.Range("C3").Value = "Cmklr1"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 1.445726333333333
$ws.Range("H3").Value = 4.337179
$ws.Range("I3").Value = 0.01544456920939864
$ws.Range("J3").Value = 0.01544456920939864
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 80.27592466666667
$ws.Range("N3").Value = 240.827774
$ws.Range("O3").Value = 0.961267782904875
$ws.Range("P3").Value = 0.9612677829048749
$ws.Range("Q3").Value = 116.0570182232829
$ws.Range("R3").Value = 1044.513164009546
$ws.Range("S3").Value = 0.01484636680183953
$ws.Range("T3").Value = 0.01484636680183953
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Rarres2"
$ws.Range("C4").Value = "Cmklr1"
$ws.Range("D4").Value = "sCs"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 1.445726333333333
$ws.Range("H4").Value = 4.337179
$ws.Range("I4").Value = 0.01544456920939864
$ws.Range("J4").Value = 0.01544456920939864
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.431192
$ws.Range("N4").Value = 1.293576
$ws.Range("O4").Value = 0.005163328601538112
$ws.Range("P4").Value = 0.005163328601538111
$ws.Range("Q4").Value = 0.6233856291226666
$ws.Range("R4").Value = 5.610470662103999
$ws.Range("S4").Value = 0.00007974538593732285
$ws.Range("T4").Value = 0.00007974538593732283
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Rarres2"
$ws.Range("C5").Value = "Cmklr1"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 75.55280566666666
$ws.Range("H5").Value = 226.658417
$ws.Range("I5").Value = 0.8071240795570661
$ws.Range("J5").Value = 0.8071240795570661
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 2.803353666666666
$ws.Range("N5").Value = 8.410060999999999
$ws.Range("O5").Value = 0.03356888849358693
$ws.Range("P5").Value = 0.03356888849358693
$ws.Range("Q5").Value = 211.801234792604
$ws.Range("R5").Value = 1906.211113133437
$ws.Range("S5").Value = 0.02709425822714014
$ws.Range("T5").Value = 0.02709425822714014
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Rarres2"
$ws.Range("C6").Value = "Cmklr1"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 75.55280566666666
$ws.Range("H6").Value = 226.658417
$ws.Range("I6").Value = 0.8071240795570661
$ws.Range("J6").Value = 0.8071240795570661
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 80.27592466666667
$ws.Range("N6").Value = 240.827774
$ws.Range("O6").Value = 0.961267782904875
$ws.Range("P6").Value = 0.9612677829048749
$ws.Range("Q6").Value = 6065.071336052639
$ws.Range("R6").Value = 54585.64202447375
$ws.Range("S6").Value = 0.7758623744849589
$ws.Range("T6").Value = 0.7758623744849588
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Rarres2"
$ws.Range("C7").Value = "Cmklr1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 75.55280566666666
$ws.Range("H7").Value = 226.658417
$ws.Range("I7").Value = 0.8071240795570661
$ws.Range("J7").Value = 0.8071240795570661
$ws.Range("K7").Value = 2
$ws.Range("L7").Value = 0.6666666666666666
$ws.Range("M7").Value = 0.431192
$ws.Range("N7").Value = 1.293576
$ws.Range("O7").Value = 0.005163328601538112
$ws.Range("P7").Value = 0.005163328601538111
$ws.Range("Q7").Value = 32.57776538102133
$ws.Range("R7").Value = 293.199888429192
$ws.Range("S7").Value = 0.004167446844967122
$ws.Range("T7").Value = 0.004167446844967121
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Rarres2"
$ws.Range("C8").Value = "Cmklr1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 16.60889166666666
$ws.Range("H8").Value = 49.82667499999999
$ws.Range("I8").Value = 0.1774313512335352
$ws.Range("J8").Value = 0.1774313512335352
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 2.803353666666666
$ws.Range("N8").Value = 8.410060999999999
$ws.Range("O8").Value = 0.03356888849358693
$ws.Range("P8").Value = 0.03356888849358693
$ws.Range("Q8").Value = 46.56059735301943
$ws.Range("R8").Value = 419.0453761771749
$ws.Range("S8").Value = 0.005956173244825001
$ws.Range("T8").Value = 0.005956173244825
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Rarres2"
$ws.Range("C9").Value = "Cmklr1"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 16.60889166666666
$ws.Range("H9").Value = 49.82667499999999
$ws.Range("I9").Value = 0.1774313512335352
$ws.Range("J9").Value = 0.1774313512335352
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 80.27592466666667
$ws.Range("N9").Value = 240.827774
$ws.Range("O9").Value = 0.961267782904875
$ws.Range("P9").Value = 0.9612677829048749
$ws.Range("Q9").Value = 1333.294136230161
$ws.Range("R9").Value = 11999.64722607145
$ws.Range("S9").Value = 0.1705590416180765
$ws.Range("T9").Value = 0.1705590416180765
$ws.Range("A10").Value = "sCs"
$ws.Range("B10").Value = "Rarres2"
$ws.Range("C10").Value = "Cmklr1"
$ws.Range("D10").Value = "sCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 16.60889166666666
$ws.Range("H10").Value = 49.82667499999999
$ws.Range("I10").Value = 0.1774313512335352
$ws.Range("J10").Value = 0.1774313512335352
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.431192
$ws.Range("N10").Value = 1.293576
$ws.Range("O10").Value = 0.005163328601538112
$ws.Range("P10").Value = 0.005163328601538111
$ws.Range("Q10").Value = 7.161621215533332
$ws.Range("R10").Value = 64.45459093979998
$ws.Range("S10").Value = 0.0009161363706336668
$ws.Range("T10").Value = 0.0009161363706336665
